$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump version from 0.1 to 1.0 (keep it stored as text, like the original "0.1")
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.0"

# The TC2/TC3/TC4 blocks get their step/result text rotated:
#   old TC2 (analisar) -> becomes "detalhar" content
#   old TC3 (cancelar) -> becomes "analisar" content
#   old TC4 (detalhar) -> becomes "cancelar" content
# A new TC3 block is effectively inserted (taking the old TC3 slot's position)
# while the former TC3 content shifts down to what is now TC4.

# Block at rows 14-18 (labelled TC2) now shows the "detalhar diaria" content
$ws.Range("B18").Value = "Beneficiário Clica em detalhar diária."
$ws.Range("D18").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"

# Block at rows 21-25 (labelled TC3) now shows the "analisar prestacao" content
$ws.Range("B25").Value = "Beneficiário Clica em analisar prestação de contas."
$ws.Range("D25").Value = "SYSTEM Apresenta a tela de Analisar Prestação de Contas"

# Block at rows 28-32 (labelled TC4) now shows the "cancelar diaria" content
$ws.Range("B32").Value = "Beneficiário Clica em cancelar diária."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"
